$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 11") {
        $shp.Delete()
    }
}
